{"js": "const replacements = [\n  [\"780\u00d78=6240\", \"676\u00d75=3380\"],\n  [\"521\u00d75=2605\", \"274\u00d76=1644\"],\n  [\"574\u00d76=3444\", \"674\u00d76=4044\"],\n  [\"403\u00d75=2015\", \"851\u00d72=1702\"],\n  [\"401\u00d79=3609\", \"908\u00d76=5448\"],\n  [\"763\u00d79=6867\", \"509\u00d72=1018\"],\n  [\"844\u00d79=7596\", \"667\u00d76=4002\"],\n  [\"480\u00d75=2400\", \"171\u00d76=1026\"],\n  [\"305\u00d79=2745\", \"338\u00d73=1014\"],\n  [\"700\u00d73=2100\", \"312\u00d73=936\"],\n  [\"762\u00d78=6096\", \"480\u00d73=1440\"],\n  [\"415\u00d77=2905\", \"623\u00d73=1869\"],\n  [\"358\u00d74=1432\", \"268\u00d79=2412\"],\n  [\"671\u00d73=2013\", \"586\u00d77=4102\"],\n  [\"572\u00d77=4004\", \"206\u00d75=1030\"],\n  [\"229\u00d73=687\", \"295\u00d76=1770\"],\n  [\"462\u00d79=4158\", \"461\u00d77=3227\"],\n  [\"114\u00d76=684\", \"272\u00d72=544\"],\n  [\"244\u00d78=1952\", \"895\u00d72=1790\"],\n  [\"769\u00d78=6152\", \"634\u00d72=1268\"],\n  [\"621\u00d79=5589\", \"846\u00d75=4230\"],\n  [\"436\u00d78=3488\", \"710\u00d73=2130\"],\n  [\"895\u00d79=8055\", \"770\u00d73=2310\"],\n  [\"982\u00d75=4910\", \"912\u00d73=2736\"],\n  [\"599\u00d76=3594\", \"725\u00d77=5075\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('780\u00d78=6240', '676\u00d75=3380'),\n    @('521\u00d75=2605', '274\u00d76=1644'),\n    @('574\u00d76=3444', '674\u00d76=4044'),\n    @('403\u00d75=2015', '851\u00d72=1702'),\n    @('401\u00d79=3609', '908\u00d76=5448'),\n    @('763\u00d79=6867', '509\u00d72=1018'),\n    @('844\u00d79=7596', '667\u00d76=4002'),\n    @('480\u00d75=2400', '171\u00d76=1026'),\n    @('305\u00d79=2745', '338\u00d73=1014'),\n    @('700\u00d73=2100', '312\u00d73=936'),\n    @('762\u00d78=6096', '480\u00d73=1440'),\n    @('415\u00d77=2905', '623\u00d73=1869'),\n    @('358\u00d74=1432', '268\u00d79=2412'),\n    @('671\u00d73=2013', '586\u00d77=4102'),\n    @('572\u00d77=4004', '206\u00d75=1030'),\n    @('229\u00d73=687', '295\u00d76=1770'),\n    @('462\u00d79=4158', '461\u00d77=3227'),\n    @('114\u00d76=684', '272\u00d72=544'),\n    @('244\u00d78=1952', '895\u00d72=1790'),\n    @('769\u00d78=6152', '634\u00d72=1268'),\n    @('621\u00d79=5589', '846\u00d75=4230'),\n    @('436\u00d78=3488', '710\u00d73=2130'),\n    @('895\u00d79=8055', '770\u00d73=2310'),\n    @('982\u00d75=4910', '912\u00d73=2736'),\n    @('599\u00d76=3594', '725\u00d77=5075'),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"No match found for: $findText\"\n    }\n}"}
